# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain overlapping event rows; the F-column counts are bumped
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for the "展览" sheet
$sheet1Updates = @{
    2  = 71
    3  = 390
    7  = 253
    8  = 13961
    9  = 101
    10 = 90
    11 = 5628
    15 = 51
    16 = 1224
    18 = 163
    19 = 756
    20 = 2907
    21 = 41
    22 = 10410
    23 = 1188
    24 = 29
    25 = 47
    26 = 3706
}

# Map of row -> new F value for the "全部类型" sheet
$sheet4Updates = @{
    2  = 71
    3  = 390
    8  = 253
    9  = 13961
    10 = 101
    11 = 90
    12 = 5628
    16 = 51
    17 = 1224
    19 = 163
    20 = 756
    21 = 2907
    22 = 41
    24 = 10410
    25 = 1188
    26 = 29
    27 = 47
    28 = 3706
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
